$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (German states' cumulative death counts) for
# 2020-05-03 (row 57) through 2020-05-17 (row 71), appended below the
# existing data which ended at row 56 (2020-05-02).
$newRows = @(
    @(57, @(43954,1412,1910,154,122,30,164,370,18,453,1287,176,138,165,44,113,93)),
    @(58, @(43955,1421,1926,154,122,30,166,372,18,456,1290,177,139,167,45,113,96)),
    @(59, @(43956,1452,1949,154,123,31,172,381,19,463,1332,182,139,171,46,118,99)),
    @(60, @(43957,1481,2001,159,127,31,190,386,19,470,1358,187,140,177,46,119,104)),
    @(61, @(43958,1497,2050,162,129,32,194,396,19,478,1372,189,140,181,48,120,109)),
    @(62, @(43959,1515,2114,163,131,32,201,403,19,493,1397,192,140,184,48,122,112)),
    @(63, @(43960,1534,2147,164,133,32,204,408,19,496,1424,195,142,186,48,122,115)),
    @(64, @(43961,1542,2153,165,134,32,204,411,19,498,1425,195,142,187,48,123,117)),
    @(65, @(43962,1545,2155,165,134,34,204,412,19,498,1437,195,142,187,48,124,118)),
    @(66, @(43963,1568,2182,165,136,35,216,415,20,507,1456,206,144,190,50,125,118)),
    @(67, @(43964,1594,2209,170,137,36,223,419,20,508,1473,207,144,191,52,125,126)),
    @(68, @(43965,1608,2229,173,138,37,228,423,20,524,1483,212,145,195,52,125,131)),
    @(69, @(43966,1628,2260,177,148,37,228,429,20,534,1493,213,147,195,54,126,135)),
    @(70, @(43967,1644,2273,181,149,37,231,432,20,535,1500,216,149,197,54,126,137)),
    @(71, @(43968,1648,2283,182,149,37,232,435,20,544,1505,216,149,196,54,126,138))
)

foreach ($entry in $newRows) {
    $rowNum = $entry[0]
    $values = $entry[1]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
    # Match the date number-format used by the existing date column (A).
    $ws.Cells.Item($rowNum, 1).NumberFormat = $ws.Cells.Item($rowNum - 1, 1).NumberFormat
}

# Move/restore the active selection to A2 (as in the updated workbook).
$null = $ws.Range("A2").Select()

# Set up the page for printing (paper size / orientation), matching the
# print setup recorded in the new workbook version.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
